$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 53; this shifts existing rows 53:84 down to 54:85
# and carries the date-style formatting on column D along with the shift.
$ws.Rows.Item(53).Insert()

# Populate the new row 53 with the new market-price record.
$ws.Cells.Item(53, 1).Value = 10
$ws.Cells.Item(53, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(53, 3).Value = "La Araucanía"
$ws.Cells.Item(53, 4).Value = 44879
$ws.Cells.Item(53, 5).Value = 9
$ws.Cells.Item(53, 6).Value = 100112026
$ws.Cells.Item(53, 7).Value = "Haba"
$ws.Cells.Item(53, 8).Value = "Sin especificar"
$ws.Cells.Item(53, 9).Value = "Primera"
$ws.Cells.Item(53, 10).Value = 170
$ws.Cells.Item(53, 11).Value = 10000
$ws.Cells.Item(53, 12).Value = 11000
$ws.Cells.Item(53, 13).Value = 10588
$ws.Cells.Item(53, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(53, 15).Value = "Región Metropolitana"
$ws.Cells.Item(53, 16).Value = 424
$ws.Cells.Item(53, 17).Value = 25
$ws.Cells.Item(53, 18).Value = "Hortaliza"
